$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A90").Value = "Sort_Drawdowns"
$ws.Range("B90").Value = "Test sort drawdown for simple returns"
$ws.Range("C90").Value = "Sort_Drawdowns_test"
